# Generate Report for Handback
#
# The "35920f53-4cbc-4f58-9b6c-676d763e3c6f" row (row 7) in both the
# zh-cn and de-de localization-status sheets receives its handback
# report: a target/handback file name (hyperlinked, like the other
# rows), the handback .xlf file name, the handback datetime, and an
# error detail noting the handback was generated from a stale handoff
# commit.

$wb = $excel.ActiveWorkbook

$fileStem = "35920f53-4cbc-4f58-9b6c-676d763e3c6f"
$currentCommit = "be26005e43e27691787d01c7db8ea445d85f3ccb"
$latestCommit  = "78ba2c552b54bdf308a6b13a2fe246c7eb676a33"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$currentCommit/e2e/$fileStem.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$latestCommit/e2e/$fileStem.md."

$sheets = @(
    @{ Name = "zh-cn"; Org = "ol-test0-zhcn"; XlfTail = "zh-cn.xlf"; HandbackTime = "2016-08-31 15:11:40" },
    @{ Name = "de-de"; Org = "ol-test0-dede"; XlfTail = "de-de.xlf"; HandbackTime = "2016-08-31 15:11:47" }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    $targetDisplay = "$fileStem.md"
    $handbackXlf = "$fileStem.a221794ca604c3686aed7871e7bafaa489be4116." + $s.XlfTail
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/" + $s.Org + "/blob/$currentCommit/e2e/$fileStem.md"

    # I7: Latest Target File -- becomes a hyperlink to the md source, like the
    # other rows in the table (A7/I2/I3/I4/I5 etc.)
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, "", "", $targetDisplay)
    $ws.Range("I7").Style = "HyperLink"

    # J7: Latest Handback File -- the xlf the handback produced
    $ws.Range("J7").Value = $handbackXlf

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $s.HandbackTime

    # P7: Error Detail -- handback came from a stale handoff commit
    $ws.Range("P7").Value = $errorDetail
}
